# Update betting-odds figures on the Flashscore "Jogos da Semana" sheet.
# Only numeric odds cells change; everything else (labels, row 1..14
# structure, styles) stays untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2 (Fortaleza vs Atletico-MG)
$ws.Range("Q2").Value = 2.25
$ws.Range("R2").Value = 1.62

# Row 5 (Sport Recife vs Operario)
$ws.Range("G5").Value = 1.57
$ws.Range("H5").Value = 3.8
$ws.Range("I5").Value = 6
$ws.Range("J5").Value = 2.25
$ws.Range("L5").Value = 6.5
$ws.Range("N5").Value = 7.5
$ws.Range("U5").Value = 2.38
$ws.Range("V5").Value = 1.53
$ws.Range("Z5").Value = 11
$ws.Range("AD5").Value = 7.5
$ws.Range("AE5").Value = 23
$ws.Range("AF5").Value = 81
$ws.Range("AH5").Value = 12
$ws.Range("AO5").Value = 8.5
$ws.Range("AQ5").Value = 29

# Row 8 (Correcaminos vs Cancun)
$ws.Range("G8").Value = 3.1
$ws.Range("H8").Value = 3.35
$ws.Range("I8").Value = 2.15
$ws.Range("J8").Value = 3.55
$ws.Range("K8").Value = 2.15
$ws.Range("L8").Value = 2.7
$ws.Range("M8").Value = 1.01
$ws.Range("N8").Value = 8.1
$ws.Range("O8").Value = 1.3
$ws.Range("P8").Value = 2.92
$ws.Range("Q8").Value = 1.95
$ws.Range("R8").Value = 1.75
$ws.Range("W8").Value = 9.25
$ws.Range("X8").Value = 16
$ws.Range("Y8").Value = 11
$ws.Range("Z8").Value = 40
$ws.Range("AA8").Value = 27
$ws.Range("AB8").Value = 35
$ws.Range("AC8").Value = 9.5
$ws.Range("AD8").Value = 6.4
$ws.Range("AE8").Value = 14.5
$ws.Range("AH8").Value = 7.4
$ws.Range("AI8").Value = 10.25
$ws.Range("AJ8").Value = 9
$ws.Range("AK8").Value = 20
$ws.Range("AL8").Value = 17.5
$ws.Range("AM8").Value = 29
$ws.Range("AN8").Value = 5
$ws.Range("AO8").Value = 16.5
$ws.Range("AP8").Value = 23
$ws.Range("AQ8").Value = 75
$ws.Range("AR8").Value = 110
$ws.Range("AT8").Value = 2.62
$ws.Range("AU8").Value = 7
$ws.Range("AW8").Value = 4.05
$ws.Range("AX8").Value = 10.75
$ws.Range("AY8").Value = 18.5
$ws.Range("AZ8").Value = 40
$ws.Range("BA8").Value = 70
$ws.Range("BB8").Value = 250

# Row 11
$ws.Range("M11").Value = 1.05
$ws.Range("O11").Value = 1.41
$ws.Range("P11").Value = 2.62

# Row 12
$ws.Range("M12").Value = 1.05
$ws.Range("O12").Value = 1.3

# Row 13
$ws.Range("M13").Value = 1.04
$ws.Range("O13").Value = 1.27

# Row 14
$ws.Range("G14").Value = 2.67
$ws.Range("I14").Value = 2.42
$ws.Range("J14").Value = 3.3
$ws.Range("L14").Value = 3
$ws.Range("O14").Value = 1.33
$ws.Range("P14").Value = 3.05
$ws.Range("Q14").Value = 2
$ws.Range("W14").Value = 8.25
$ws.Range("X14").Value = 13
$ws.Range("Z14").Value = 30
$ws.Range("AH14").Value = 8
$ws.Range("AI14").Value = 12
$ws.Range("AJ14").Value = 9.5
$ws.Range("AK14").Value = 26
$ws.Range("AM14").Value = 30
$ws.Range("AN14").Value = 4.65
$ws.Range("AO14").Value = 14.5
$ws.Range("AQ14").Value = 65
$ws.Range("AW14").Value = 4.35
$ws.Range("AX14").Value = 13
$ws.Range("AZ14").Value = 55
